# Adds a new "Output mong muốn" worksheet right after "Cover " (becoming the
# 2nd tab), populates it with the checklist-authoring guidance content, and
# makes it the active tab - matching the author's commit ("#feat: add rule file").

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new worksheet right after "Cover " -----------------------
$coverSheet = $wb.Worksheets.Item("Cover ")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $coverSheet)
$ws.Name = "Output mong muốn"

# Turn off gridlines to match the rest of the workbook's sheets look & feel.
$excel.ActiveWindow.DisplayGridlines = $false

# --- 2. Section headers (bold, light-gray fill) ------------------------------
$headerCells = @("A2", "A7", "A13", "A20")
foreach ($addr in $headerCells) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Interior.Color = 13948888
}

$ws.Range("A2").Value = "1. Trước khi đưa code lên"
$ws.Range("A7").Value = "2. Phần coding"
$ws.Range("A13").Value = "3. Phần Unit test"
$ws.Range("A20").Value = "4. Hướng dẫn"

# --- 3. Bullet sub-items under each section ----------------------------------
$ws.Range("B3").Value = "- Check hết tất cả rule trong checklist"
$ws.Range("B4").Value = "- Chạy lint để không còn warning"
$ws.Range("B5").Value = "- Cài eslint và prettier format code"

$ws.Range("B8").Value = "- Sau khi code xong phải tối ưu lại code nếu có"
$ws.Range("B9").Value = "- Đặt tên biến rõ ràng đúng ý nghĩa, add tsdoc rõ ràng"

$ws.Range("B14").Value = "- Yêu cầu coverage >= 80% "
$ws.Range("B15").Value = "- Code rõ ràng, xóa hết code thừa k sử dụng"

$ws.Range("A22").Value = "- Project sử dụng lib primevue, primeflex => Tận dụng hết tất cả component của lib trước"
$ws.Range("A23").Value = "- Project đang viết vue3 setup script nên viết theo đúng kiểu của code mẫu đang sử dụng"
$ws.Range("A24").Value = "- Unit test sử dụng mocha + sinon + chai => có gì thì search theo các keyword này để tìm hiểu"
$ws.Range("A25").Value = "- Tham khảo cấu trúc thư mục hiện tại để làm theo"

$bulletCells = @("B3","B4","B5","B8","B9","B14","B15","A22","A23","A24","A25")
foreach ($addr in $bulletCells) {
    $ws.Range($addr).HorizontalAlignment = -4131
}

# --- 4. Column widths (match the other sheets' default 9.14 width) ----------
$ws.Columns.ColumnWidth = 9.140625

# --- 5. Make the new sheet the active / visible tab --------------------------
$ws.Activate()
$wb.Windows.Item(1).ActiveSheet.Range("A2").Select()
